$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Copy the formatting from column J (the existing "LfoAmp" parameter column)
# onto the new column K so the new "Smackidiboo" column matches the look
# of the other parameter columns (Neutral style, borders, alignment, etc.).
[void]$ws.Range("J2:J9").Copy()
[void]$ws.Range("K2:K9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New header for the added output parameter column.
$ws.Range("K2").Value = "Smackidiboo"

# New parameter values, mirroring the other 0-initialized parameter rows.
$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("K9").Value = 0

# Update the selection to reflect where the user was working in the UI.
[void]$ws.Range("K3").Select()
